# Update "paises" (countries) dashboard: refresh timestamp, swap Chile/España
# and Kuwait/Paises Bajos rank order, and refresh case-count figures for the
# countries whose totals moved since the last snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "as of" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 18:24"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 3003464
$ws.Range("C4").Value = 20536
$ws.Range("D4").Value = 1292077
$ws.Range("E4").Value = 1578709
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 109
$ws.Range("H4").Value = 132678

# --- Row 6: India ---
$ws.Range("B6").Value = 715190
$ws.Range("C6").Value = 17354
$ws.Range("D6").Value = 437189
$ws.Range("E6").Value = 257871
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 430
$ws.Range("H6").Value = 20130

# --- Row 9: España -> Chile overtakes España in the ranking ---
$ws.Range("A9").Value = "Chile"
$ws.Range("B9").Value = 298557
$ws.Range("C9").Value = 3025
$ws.Range("D9").Value = 264371
$ws.Range("E9").Value = 27802
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 76
$ws.Range("H9").Value = 6384

# --- Row 10: Chile -> España drops to this row ---
$ws.Range("A10").Value = "España"
$ws.Range("B10").Value = 297625
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 28385

# --- Row 11: Reino Unido ---
$ws.Range("B11").Value = 285768
$ws.Range("C11").Value = 352
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 16
$ws.Range("H11").Value = 44236

# --- Row 14: Italia ---
$ws.Range("B14").Value = 241819
$ws.Range("C14").Value = 208
$ws.Range("D14").Value = 192241
$ws.Range("E14").Value = 14709
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 34869

# --- Row 23: Canada ---
$ws.Range("B23").Value = 105764
$ws.Range("C23").Value = 228
$ws.Range("D23").Value = 69431
$ws.Range("E23").Value = 27646
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 8687

# --- Row 35: Kuwait -> Paises Bajos overtakes Kuwait in the ranking ---
$ws.Range("A35").Value = "Paises Bajos"
$ws.Range("B35").Value = 50657
$ws.Range("C35").Value = 36
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 6128

# --- Row 36: Paises Bajos -> Kuwait drops to this row ---
$ws.Range("A36").Value = "Kuwait"
$ws.Range("B36").Value = 50644
$ws.Range("C36").Value = 703
$ws.Range("D36").Value = 41001
$ws.Range("E36").Value = 9270
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 373

# --- Row 41: Singapur ---
$ws.Range("D41").Value = 40717
$ws.Range("E41").Value = 4240

# --- Row 61: Moldavia ---
$ws.Range("B61").Value = 17906
$ws.Range("C61").Value = 92
$ws.Range("E61").Value = 6267
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 592

# --- Row 69: Chequia ---
$ws.Range("B69").Value = 12532
$ws.Range("C69").Value = 17
$ws.Range("D69").Value = 7873
$ws.Range("E69").Value = 4309
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 350

# --- Row 110: Cuba ---
$ws.Range("B110").Value = 2380
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 2234
$ws.Range("E110").Value = 60

# --- Row 111: Mali ---
$ws.Range("B111").Value = 2331
$ws.Range("C111").Value = 1
$ws.Range("D111").Value = 1547
$ws.Range("E111").Value = 665

# --- Row 131: Jordania ---
$ws.Range("B131").Value = 1167
$ws.Range("C131").Value = 3
$ws.Range("D131").Value = 957
$ws.Range("E131").Value = 200
